$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 : GA02 / Gabriel Simard / 15-09-2015 / NomJeu de varchar(30) à varchar(50) / Terminée
$ws.Range("A10").Value = "GA02"
$ws.Range("B10").Value = "Gabriel Simard"
$ws.Range("C10").NumberFormat = "m/d/yy"
$ws.Range("C10").Value2 = 42262
$ws.Range("D10").Value = "NomJeu de varchar(30) à varchar(50)"
$ws.Range("E10").Value = "Terminée"

# Row 11 : EL01 / Élodie Kérouak / 15-09-2015 / Tag devient null / Terminée
$ws.Range("A11").Value = "EL01"
$ws.Range("B11").Value = "Élodie Kérouak"
$ws.Range("C11").NumberFormat = "m/d/yy"
$ws.Range("C11").Value2 = 42262
$ws.Range("D11").Value = "Tag devient null"
$ws.Range("E11").Value = "Terminée"

# Update the selected cell shown in the saved sheet view
[void]$ws.Range("G9").Select()
